$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transacciones")

# ---------------------------------------------------------------------------
# Row 89: "Mes de Gimnasio" / "Mara Sport" (Cuidado Personal, Tarjeta Banamex)
# ---------------------------------------------------------------------------
$ws.Cells.Item(89, 1).Value = 43577
$ws.Cells.Item(89, 2).Value = 500
$ws.Cells.Item(89, 3).Value = "Mes de Gimnasio"
$ws.Cells.Item(89, 4).Value = "Cuidado Personal"
$ws.Cells.Item(89, 5).Value = "Gasto"
$ws.Cells.Item(89, 6).Value = "Tarjeta Banamex"
$ws.Cells.Item(89, 7).Value = "Mara Sport"
$ws.Cells.Item(89, 11).Formula = "=K88-B89"
$ws.Cells.Item(89, 12).Value = 2442.5700000000002
$ws.Cells.Item(89, 13).Value = 2
$ws.Cells.Item(89, 14).Formula = "=SUM(K89:M89)"
$ws.Cells.Item(89, 15).Formula = "=N89-4000"
$ws.Cells.Item(89, 16).Formula = "=O89-Ahorros!`$E`$4"

# Match the date + highlight styles used on the row directly above it.
$ws.Cells.Item(88, 1).Copy()
$ws.Cells.Item(89, 1).PasteSpecial(-4122)
$ws.Cells.Item(88, 16).Copy()
$ws.Cells.Item(89, 16).PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Row 90: "Coca Cola " / "Oxxo" (Golosinas, Tarjeta Banamex)
# ---------------------------------------------------------------------------
$ws.Cells.Item(90, 1).Value = 43577
$ws.Cells.Item(90, 2).Value = 13.5
$ws.Cells.Item(90, 3).Value = "Coca Cola "
$ws.Cells.Item(90, 4).Value = "Golosinas"
$ws.Cells.Item(90, 5).Value = "Gasto"
$ws.Cells.Item(90, 6).Value = "Tarjeta Banamex"
$ws.Cells.Item(90, 7).Value = "Oxxo"
$ws.Cells.Item(90, 11).Formula = "=K89-B90"
$ws.Cells.Item(90, 12).Value = 2442.5700000000002
$ws.Cells.Item(90, 13).Value = 2
$ws.Cells.Item(90, 14).Formula = "=SUM(K90:M90)"
$ws.Cells.Item(90, 15).Formula = "=N90-4000"
$ws.Cells.Item(90, 16).Formula = "=O90-Ahorros!`$E`$4"

$ws.Cells.Item(88, 1).Copy()
$ws.Cells.Item(90, 1).PasteSpecial(-4122)
$ws.Cells.Item(88, 16).Copy()
$ws.Cells.Item(90, 16).PasteSpecial(-4122)
$ws.Cells.Item(89, 14).Copy()
$ws.Cells.Item(90, 14).PasteSpecial(-4122)
$ws.Cells.Item(89, 15).Copy()
$ws.Cells.Item(90, 15).PasteSpecial(-4122)

# Move the active selection to Q90, matching the sheet's last interacted cell.
$ws.Range("Q90").Select() | Out-Null

$excel.CutCopyMode = $false
